$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.687.77'
$ws.Range("E2").Value = '  +3.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.013.85'
$ws.Range("E3").Value = '  +2.23%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '528.01'
$ws.Range("E5").Value = '  +2.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.59'
$ws.Range("E6").Value = '  +1.86%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  +1.18%  '

$ws.Range("E8").Value = '  +0.22%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.739'
$ws.Range("E9").Value = '  +1.50%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.177'
$ws.Range("E10").Value = '  +2.94%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000346'
$ws.Range("E11").Value = '  +1.15%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '44.61'
$ws.Range("E12").Value = '  +3.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.77'
$ws.Range("E13").Value = '  +3.85%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.642.52'
$ws.Range("E14").Value = '  +1.75%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.981.54'
$ws.Range("E15").Value = '  +1.44%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.39'
$ws.Range("E16").Value = '  +8.61%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.32'
$ws.Range("E17").Value = '  +1.76%  '

$ws.Range("E18").Value = '  +0.77%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.134'
$ws.Range("E19").Value = '  -1.65%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.581.79'
$ws.Range("E20").Value = '  +3.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '441.47'
$ws.Range("E21").Value = '  +2.31%  '

$ws.Range("E22").Value = '  +5.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '94.05'
$ws.Range("E23").Value = '  +6.80%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.39'
$ws.Range("E24").Value = '  -0.55%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.29'
$ws.Range("E25").Value = '  +4.46%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.11'
$ws.Range("E26").Value = '  +5.95%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.93'
$ws.Range("E27").Value = '  -0.92%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37.23'
$ws.Range("E28").Value = '  +1.47%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '704.92'
$ws.Range("E29").Value = '  +0.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.63'
$ws.Range("E30").Value = '  +2.79%  '

$ws.Range("E31").Value = '  +1.82%  '

$ws.Range("E32").Value = '  +1.08%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.99'
$ws.Range("E33").Value = '  +17.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '68.54'
$ws.Range("E34").Value = '  +3.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0914'
$ws.Range("E35").Value = '  +5.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.448'
$ws.Range("E36").Value = '  +1.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '41.04'
$ws.Range("E37").Value = '  +1.98%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.67'
$ws.Range("E38").Value = '  +21.53%  '

$ws.Range("E39").Value = '  +2.00%  '

$ws.Range("E40").Value = '  -0.17%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0497'
$ws.Range("E41").Value = '  +2.89%  '

$ws.Range("E42").Value = '  -0.20%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.90'
$ws.Range("E43").Value = '  +2.96%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.15'
$ws.Range("E44").Value = '  +0.82%  '

$ws.Range("E45").Value = '  +6.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.23'
$ws.Range("E46").Value = '  +9.72%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.147'
$ws.Range("E47").Value = '  +2.76%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.000284'
$ws.Range("E48").Value = '  +22.50%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.33'
$ws.Range("E49").Value = '  +7.27%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.39'
$ws.Range("E50").Value = '  +1.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0343'
$ws.Range("E51").Value = '  -3.38%  '
